$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("J2").Value = 6783
$ws.Range("J3").Value = 7174
$ws.Range("I4").Value = 1773
$ws.Range("J4").Value = 1563
$ws.Range("J5").Value = 561
$ws.Range("J6").Value = 9561
$ws.Range("I7").Value = 26231
$ws.Range("J7").Value = 25642

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J3").Value = 45
$ws.Range("J6").Value = 243
$ws.Range("J7").Value = 368

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 482
$ws.Range("J4").Value = 86
$ws.Range("J6").Value = 568
$ws.Range("J7").Value = 1611

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("J2").Value = 152
$ws.Range("J7").Value = 512

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 269
$ws.Range("J3").Value = 386
$ws.Range("J6").Value = 403
$ws.Range("J7").Value = 1159

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("J3").Value = 133
$ws.Range("J7").Value = 371

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 186
$ws.Range("J3").Value = 183
$ws.Range("J6").Value = 231
$ws.Range("J7").Value = 643

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("J5").Value = 77
$ws.Range("J6").Value = 195
$ws.Range("J7").Value = 743
$ws.Range("J8").Value = 1611
$ws.Range("J9").Value = 137
$ws.Range("J10").Value = 187
$ws.Range("J15").Value = 304
$ws.Range("J18").Value = 216
$ws.Range("J19").Value = 750
$ws.Range("J20").Value = 538
$ws.Range("J22").Value = 60
$ws.Range("J26").Value = 52
$ws.Range("J27").Value = 151
$ws.Range("J29").Value = 1390
$ws.Range("J33").Value = 1159
$ws.Range("J36").Value = 350
$ws.Range("J39").Value = 18
$ws.Range("J42").Value = 1104
$ws.Range("J44").Value = 196
$ws.Range("J47").Value = 190
$ws.Range("J49").Value = 162
$ws.Range("J51").Value = 312
$ws.Range("J53").Value = 368
$ws.Range("J54").Value = 492
$ws.Range("J58").Value = 16
$ws.Range("J63").Value = 77
$ws.Range("J65").Value = 643
$ws.Range("J66").Value = 79
$ws.Range("J67").Value = 964
$ws.Range("J73").Value = 246
$ws.Range("J75").Value = 76
$ws.Range("J78").Value = 304
$ws.Range("J83").Value = 512
$ws.Range("J85").Value = 1059
$ws.Range("I86").Value = 168
$ws.Range("J86").Value = 164
$ws.Range("J88").Value = 267
$ws.Range("J89").Value = 328
$ws.Range("J91").Value = 296
$ws.Range("J94").Value = 275
$ws.Range("J95").Value = 371
$ws.Range("J97").Value = 235
$ws.Range("J98").Value = 188
$ws.Range("I101").Value = 26231
$ws.Range("J101").Value = 25642

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J2").Value = 245
$ws.Range("J3").Value = 356
$ws.Range("J4").Value = 67
$ws.Range("J6").Value = 269
$ws.Range("J7").Value = 964

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("J6").Value = 92
$ws.Range("J7").Value = 162

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("J6").Value = 230
$ws.Range("J7").Value = 492

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("J2").Value = 421
$ws.Range("J3").Value = 493
$ws.Range("J7").Value = 1390

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J3").Value = 215
$ws.Range("J6").Value = 290
$ws.Range("J7").Value = 750

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("J6").Value = 77
$ws.Range("J7").Value = 196

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J2").Value = 58
$ws.Range("J7").Value = 195

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 219
$ws.Range("J7").Value = 1104

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("J6").Value = 106
$ws.Range("J7").Value = 187

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J2").Value = 79
$ws.Range("J4").Value = 32
$ws.Range("J5").Value = 5
$ws.Range("J6").Value = 94
$ws.Range("J7").Value = 304

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("J2").Value = 78
$ws.Range("J4").Value = 11
$ws.Range("J7").Value = 296

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("J2").Value = 150
$ws.Range("J7").Value = 538

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J3").Value = 44
$ws.Range("J7").Value = 216

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J3").Value = 112
$ws.Range("J7").Value = 350

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 229
$ws.Range("J3").Value = 227
$ws.Range("J7").Value = 743

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("J6").Value = 148
$ws.Range("J7").Value = 275

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J2").Value = 43
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 190

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 85
$ws.Range("J7").Value = 304

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J6").Value = 119
$ws.Range("J7").Value = 188

$ws = $wb.Worksheets.Item('East Village')
$ws.Range("J2").Value = 11
$ws.Range("J7").Value = 52

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("J3").Value = 6
$ws.Range("J6").Value = 18

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J2").Value = 15
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J2").Value = 36
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 81
$ws.Range("J7").Value = 246

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J4").Value = 8
$ws.Range("J6").Value = 162
$ws.Range("J7").Value = 235

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J3").Value = 64
$ws.Range("J6").Value = 134
$ws.Range("J7").Value = 267

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J4").Value = 34
$ws.Range("J7").Value = 328

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J3").Value = 35
$ws.Range("J7").Value = 151

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J3").Value = 28
$ws.Range("I4").Value = 79
$ws.Range("I7").Value = 168
$ws.Range("J7").Value = 164

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("J2").Value = 32
$ws.Range("J3").Value = 25
$ws.Range("J7").Value = 76

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 71
$ws.Range("J7").Value = 312

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J3").Value = 371
$ws.Range("J6").Value = 308
$ws.Range("J7").Value = 1059

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("J2").Value = 29
$ws.Range("J7").Value = 60

$ws = $wb.Worksheets.Item('Millenium Park')
$ws.Range("J2").Value = 2
$ws.Range("J7").Value = 16
